# Apply text corrections described in commit "From 1.2.4 to 1.2.5 change and minor updates"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Version bump: 0.1 -> 1.2.5
$ws.Range("D2").Value = "1.2.5"

# 2. Precondition text fix (appears once per test case, same shared text)
$oldPrecond = "O usuario devidamente autenticado e na tela de listagem de empenhos"
$newPrecond = "O usuário devidamente autenticado e na tela de listagem de empenhos."
$precondCells = @("B8", "B16", "B24", "B31", "B38", "B45", "B53")
foreach ($cell in $precondCells) {
    $ws.Range($cell).Value = $newPrecond
}

# 3. Step text fix: capitalize "Filtra" and add trailing period (shared text, multiple cells)
$newFiltra = "Chefe/Beneficiário Filtra a listagem por registros cujos beneficiários não possuem número do credor."
$filtraCells = @("B10", "B47", "B55")
foreach ($cell in $filtraCells) {
    $ws.Range($cell).Value = $newFiltra
}

# 4. TC2 expected result text update
$ws.Range("D18").Value = "SYSTEM Exibe a lista de solicitações aguardando serem empenhadas, de todos os servidores, ordenado pelo número da diária em ordem crescente."

# 5. TC3 expected result - add trailing period
$ws.Range("D26").Value = "SYSTEM Apresenta a tela de Registrar Empenho."

# 6. TC5 expected result - add trailing period
$ws.Range("D40").Value = "SYSTEM Recupera e exibe todos os detalhes (dados) da solicitação para o usuário; e Apresenta a tela de Detalhar Diárias."
